$d = $word.ActiveDocument

# ============================================================
# Process document bottom-to-top to keep paragraph indices stable.
# ============================================================

# --- 1. Remove trailing empty paragraph before "Personal Details" (para 98) ---
$d.Paragraphs(98).Range.Delete()

# --- 2. Remove trailing empty paragraph before "Certification" (para 95) ---
$d.Paragraphs(95).Range.Delete()

# --- 3. Remove trailing empty paragraph before "Achievements" (para 84) ---
$d.Paragraphs(84).Range.Delete()

Write-Host "After step 1-3, paragraph count: " $d.Paragraphs.Count
Write-Host "Para 83: [" $d.Paragraphs(83).Range.Text "]"
Write-Host "Para 84: [" $d.Paragraphs(84).Range.Text "]"

# --- 4. Role as DevOps Engineer bullets (74-81): tense fixes ---
$d.Content.Find.Execute("AWS Infrastructure Development: Building robust AWS infrastructure from scratch using Terraform, improving deployment efficiency and scalability. Automating provisioning of VPC, EC2, ALB, EBS, NLB, IAM, and S3 resources.", $true, $false, $false, $false, $false, $true, 1, $false, "AWS Infrastructure Development: Built robust AWS infrastructure from scratch using Terraform, improving deployment efficiency and scalability. Automated provisioning of VPC, EC2, ALB, EBS, NLB, IAM, and S3 resources.", 2) | Out-Null

$d.Content.Find.Execute("Microservices Containerization: Implementing Docker-based containerization and orchestrated services with Kubernetes, ensuring high availability. Creating and managing Docker images and containers for consistent deployment.", $true, $false, $false, $false, $false, $true, 1, $false, "Microservices Containerization: Implemented Docker-based containerization and orchestrated services with Kubernetes, ensuring high availability. Created and managed Docker images and containers for consistent deployment.", 2) | Out-Null

$d.Content.Find.Execute("Configuration Management: Utilizing Ansible for consistent setup across environments. Developing and maintaining Docker files, aiding in creating base images for QA and Performance teams.", $true, $false, $false, $false, $false, $true, 1, $false, "Configuration Management: Utilized Ansible for consistent setup across environments. Developed and maintained Docker files, aiding in creating base images for QA and Performance teams.", 2) | Out-Null

$d.Content.Find.Execute("Monitoring and Observability: Setting up monitoring with Prometheus and Grafana, enhancing visibility into system health. Creating scripts for proactive issue resolution across build machines, CI/CD applications, and Kubernetes clusters.", $true, $false, $false, $false, $false, $true, 1, $false, "Monitoring and Observability: Set up monitoring with Prometheus and Grafana, enhancing visibility into system health. Created scripts for proactive issue resolution across build machines, CI/CD applications, and Kubernetes clusters.", 2) | Out-Null

$d.Content.Find.Execute("CI/CD Automation: Preparing end-to-end pipelines with Jenkins, automating software releases and minimizing downtime. Managing Jenkins pipelines for integration, nightly, and release builds, reducing manual effort.", $true, $false, $false, $false, $false, $true, 1, $false, "CI/CD Automation: Prepared end-to-end pipelines with Jenkins, automating software releases and minimizing downtime. Managed Jenkins pipelines for integration, nightly, and release builds, reducing manual effort.", 2) | Out-Null

$d.Content.Find.Execute("Version Control Management: Managing GitLab repositories, managed merge requests, and creating branches according to project requirements. Utilizing GitHub for collaboration and version control.", $true, $false, $false, $false, $false, $true, 1, $false, "Version Control Management: Managed GitLab repositories, managed merge requests, and created branches according to project requirements. Utilized GitHub for collaboration and version control.", 2) | Out-Null

$d.Content.Find.Execute("Security and Code Quality: Integrating SonarQube for code analysis, ensuring code quality and security. Conducting Trivy scans for Docker images to identify vulnerabilities early in the development lifecycle.", $true, $false, $false, $false, $false, $true, 1, $false, "Security and Code Quality: Integrated SonarQube for code analysis, ensuring code quality and security. Conducted Trivy scans for Docker images to identify vulnerabilities early in the development lifecycle.", 2) | Out-Null

$d.Content.Find.Execute("Kubernetes Expertise: Working with Kubernetes resources including deployments, namespaces, persistent volumes, services, and configmaps. Developing YAML configurations for infrastructure as code management.", $true, $false, $false, $false, $false, $true, 1, $false, "Kubernetes Expertise: Worked with Kubernetes resources including deployments, namespaces, persistent volumes, services, and configmaps. Developed YAML configurations for infrastructure as code management.", 2) | Out-Null

Write-Host "Step 4 done. Para 74: [" $d.Paragraphs(74).Range.Text "]"

# --- 5. Remove empty paragraph (72) between Cross-functional Collaboration bullet and "Role as DevOps Engineer:" ---
$d.Paragraphs(72).Range.Delete()

# --- 6. Role as Cloud Engineer bullets (61-71): tense fixes ---
$d.Content.Find.Execute("AWS Infrastructure Management: Provisioning and managing a comprehensive suite of AWS services (EC2, IAM, S3, RDS, Lambda, VPC, Route 53, EKS) to support scalable application deployment. Optimizing storage solutions by configuring EBS volumes and utilizing AWS LVM for efficient volume management.", $true, $false, $false, $false, $false, $true, 1, $false, "AWS Infrastructure Management: Provisioned and managed a comprehensive suite of AWS services (EC2, IAM, S3, RDS, Lambda, VPC, Route 53, EKS) to support scalable application deployment. Optimized storage solutions by configuring EBS volumes and utilizing AWS LVM for efficient volume management.", 2) | Out-Null

$d.Content.Find.Execute("Containerization and Orchestration: Spearheading containerization using Docker and orchestrating microservices with Kubernetes, ensuring high availability and operational efficiency. Developing and maintaining Docker images and containers for consistent deployment practices.", $true, $false, $false, $false, $false, $true, 1, $false, "Containerization and Orchestration: Spearheaded containerization using Docker and orchestrated microservices with Kubernetes, ensuring high availability and operational efficiency. Developed and maintained Docker images and containers for consistent deployment practices.", 2) | Out-Null

$d.Content.Find.Execute("CI/CD Pipeline Design: Designing and maintaining CI/CD pipelines with Jenkins, enhancing software delivery lifecycle efficiency. Collaborating with development teams to establish Jenkins pipelines, achieving faster release cycles.", $true, $false, $false, $false, $false, $true, 1, $false, "CI/CD Pipeline Design: Designed and maintained CI/CD pipelines with Jenkins, enhancing software delivery lifecycle efficiency. Collaborated with development teams to establish Jenkins pipelines, achieving faster release cycles.", 2) | Out-Null

$d.Content.Find.Execute("Security and VPN Configuration: Setting up secure communication channels through Site-to-Site and Point-to-Site VPN connections. Collaborating with security teams to implement IAM roles and policies, enhancing data protection and compliance.", $true, $false, $false, $false, $false, $true, 1, $false, "Security and VPN Configuration: Set up secure communication channels through Site-to-Site and Point-to-Site VPN connections. Collaborated with security teams to implement IAM roles and policies, enhancing data protection and compliance.", 2) | Out-Null

$d.Content.Find.Execute("Infrastructure as Code (IaC): Developing AWS infrastructure using Terraform, streamlining deployment and enhancing scalability through reusable modules.", $true, $false, $false, $false, $false, $true, 1, $false, "Infrastructure as Code (IaC): Developed AWS infrastructure using Terraform, streamlining deployment and enhancing scalability through reusable modules.", 2) | Out-Null

$d.Content.Find.Execute("Version Control and Collaboration: Utilizing Git, GitLab, and GitHub for continuous development workflow. Managing branches for Docker files, Ansible playbooks, and Kubernetes manifests based on project requirements.", $true, $false, $false, $false, $false, $true, 1, $false, "Version Control and Collaboration: Utilized Git, GitLab, and GitHub for continuous development workflow. Managed branches for Docker files, Ansible playbooks, and Kubernetes manifests based on project requirements.", 2) | Out-Null

$d.Content.Find.Execute("Linux File System Management: Creating file systems with LVM, developing physical volumes (PVs) and volume groups (VGs) tailored to project specifications.", $true, $false, $false, $false, $false, $true, 1, $false, "Linux File System Management: Created file systems with LVM, developing physical volumes (PVs) and volume groups (VGs) tailored to project specifications.", 2) | Out-Null

$d.Content.Find.Execute("System Monitoring and Maintenance: Implementing AWS Systems Manager for automated patch management and use AWS CloudWatch to monitor resources, proactively addressing performance issues.", $true, $false, $false, $false, $false, $true, 1, $false, "System Monitoring and Maintenance: Implemented AWS Systems Manager for automated patch management and used AWS CloudWatch to monitor resources, proactively addressing performance issues.", 2) | Out-Null

$d.Content.Find.Execute("Backup and Recovery: Establishing and managing backup protocols using AWS Backup and snapshots, ensuring data integrity and availability.", $true, $false, $false, $false, $false, $true, 1, $false, "Backup and Recovery: Established and managed backup protocols using AWS Backup and snapshots, ensuring data integrity and availability.", 2) | Out-Null

$d.Content.Find.Execute("Cost Optimization and Resource Utilization: Analyzing AWS resource usage, leading cost-reduction efforts through automation and best practices for resource tagging and rightsizing.", $true, $false, $false, $false, $false, $true, 1, $false, "Cost Optimization and Resource Utilization: Analyzed AWS resource usage, leading cost-reduction efforts through automation and best practices for resource tagging and rightsizing.", 2) | Out-Null

$d.Content.Find.Execute("Cross-functional Collaboration: Engaging actively with development and operations teams to align on project objectives and delivery timelines, driving continuous improvement initiatives.", $true, $false, $false, $false, $false, $true, 1, $false, "Cross-functional Collaboration: Engaged actively with development and operations teams to align on project objectives and delivery timelines, driving continuous improvement initiatives.", 2) | Out-Null

Write-Host "Step 6 done. Para 61: [" $d.Paragraphs(61).Range.Text "]"
Write-Host "Para 71: [" $d.Paragraphs(71).Range.Text "]"
Write-Host "Para 72: [" $d.Paragraphs(72).Range.Text "]"

# --- 7. Senior Infra Developer block (paras 42-53) ---
# 42: remove space "Cognizant, Noida" -> "Cognizant,Noida"
$d.Content.Find.Execute("Since Dec’24 with Cognizant, Noida", $true, $false, $false, $false, $false, $true, 1, $false, "Since Dec’24 with Cognizant,Noida", 2) | Out-Null

# 44: shorten to "Role as Senior Infra Developer:"; 45 (currently empty) becomes the new ListBullet line
$p44 = $d.Paragraphs(44)
$p44.Range.Text = "Role as Senior Infra Developer:"
$p45 = $d.Paragraphs(45)
$p45.Style = "List Bullet"
$p45.Range.Text = "Project Involvement – Multi-Region Migration & DevOps Enablement"

# 46-53: tense fixes
$d.Content.Find.Execute("Migration & Regional Expansion: Contributing to the migration and infrastructure expansion across multiple AWS regions. Actively involved in provisioning new environments to support regional deployments.", $true, $false, $false, $false, $false, $true, 1, $false, "Migration & Regional Expansion: Contributed to the migration and infrastructure expansion across multiple AWS regions. Actively involved in provisioning new environments to support regional deployments.", 2) | Out-Null

$d.Content.Find.Execute("Helm Upgrade & Management: Upgrading Helm versions for new regions, ensuring compatibility with application deployment pipelines and cluster configurations.", $true, $false, $false, $false, $false, $true, 1, $false, "Helm Upgrade & Management: Upgraded Helm versions for new regions, ensuring compatibility with application deployment pipelines and cluster configurations.", 2) | Out-Null

$d.Content.Find.Execute("CI/CD Enhancements: Updating Jenkinsfiles to support multi-region and multi-environment builds using dynamic parameters and variables, streamlining deployment processes.", $true, $false, $false, $false, $false, $true, 1, $false, "CI/CD Enhancements: Updated Jenkinsfiles to support multi-region and multi-environment builds using dynamic parameters and variables, streamlining deployment processes.", 2) | Out-Null

$d.Content.Find.Execute("Secret Management: Managing and updating secrets for newly added regions to maintain secure access and compliance.", $true, $false, $false, $false, $false, $true, 1, $false, "Secret Management: Managed and updated secrets for newly added regions to maintain secure access and compliance.", 2) | Out-Null

$d.Content.Find.Execute("Infrastructure Provisioning with Terraform: Creating and managing Kubernetes clusters in new regions using Terraform, ensuring consistency and scalability through reusable modules.", $true, $false, $false, $false, $false, $true, 1, $false, "Infrastructure Provisioning with Terraform: Created and managed Kubernetes clusters in new regions using Terraform, ensuring consistency and scalability through reusable modules.", 2) | Out-Null

$d.Content.Find.Execute("Monitoring & Observability: Implementing and maintaining monitoring solutions using AWS CloudWatch and Datadog to ensure real-time visibility into infrastructure and application performance.", $true, $false, $false, $false, $false, $true, 1, $false, "Monitoring & Observability: Implemented and maintained monitoring solutions using AWS CloudWatch and Datadog to ensure real-time visibility into infrastructure and application performance.", 2) | Out-Null

$d.Content.Find.Execute("Version Control & Collaboration: Creating feature branches in GitHub as per project requirements and collaborating with team members to review and merge changes into the main branch following standard Git practices.", $true, $false, $false, $false, $false, $true, 1, $false, "Version Control & Collaboration: Created feature branches in GitHub as per project requirements and collaborated with team members to review and merge changes into the main branch following standard Git practices.", 2) | Out-Null

$d.Content.Find.Execute("Terraform Code Enhancements: Refactoring Terraform configurations to support multi-region deployments and automate provisioning across environments.", $true, $false, $false, $false, $false, $true, 1, $false, "Terraform Code Enhancements: Refactored Terraform configurations to support multi-region deployments and automate provisioning across environments.", 2) | Out-Null

Write-Host "Step 7 done."
for ($i = 42; $i -le 53; $i++) {
    Write-Host $i ": [" $d.Paragraphs($i).Range.Text "]"
}

# --- 8. Technical Skills section (paras 35-38): merge 4 lines into 2 ---
$p35 = $d.Paragraphs(35)
$p35.Range.Text = "Infrastructure as Code: Terraform Configuration Management: Ansible"
$d.Paragraphs(36).Range.Delete()
$p36 = $d.Paragraphs(36)
$p36.Range.Text = "Programming Languages: Python, Bash Operating Systems: Windows, Linux/Ubuntu"
$d.Paragraphs(37).Range.Delete()

Write-Host "Step 8 done."
Write-Host "35: [" $d.Paragraphs(35).Range.Text "]"
Write-Host "36: [" $d.Paragraphs(36).Range.Text "]"
Write-Host "37: [" $d.Paragraphs(37).Range.Text "]"

# --- 9. Core Competencies section (paras 10-27): 18 plain paragraphs -> 11 ListBullet paragraphs ---
$ccTexts = @(
  "Cloud Infrastructure Management",
  "Configuration Management",
  "CI/CD Pipeline Automation",
  "Change Management",
  "Infrastructure as Code (IaC)",
  "Cloud Security & Compliance",
  "Containerization and Orchestration",
  "Risk Assessment",
  "ELK Stack: Implemented and maintained ELK stack for centralized logging and monitoring, enabling efficient debugging and performance analysis.",
  "Prometheus: Implemented and maintained Prometheus for application and infrastructure monitoring.",
  "Grafana: Implemented and maintained Grafana dashboards for visualizing metrics and logs."
)
$ccStart = 10
for ($i = 0; $i -lt $ccTexts.Length; $i++) {
    $p = $d.Paragraphs($ccStart + $i)
    $p.Style = "List Bullet"
    $p.Range.Text = $ccTexts[$i]
}
# Original block had 18 paragraphs (10-27); 11 are now reused (10-20); remove the remaining 7 (now sitting at position 21 repeatedly)
for ($i = 0; $i -lt 7; $i++) {
    $d.Paragraphs(21).Range.Delete()
}

Write-Host "Step 9 done. Count: " $d.Paragraphs.Count
for ($i = 8; $i -le 22; $i++) {
    Write-Host $i ": [" $d.Paragraphs($i).Range.Text "]"
}

# --- 10. Profile Summary paragraph (para 6) -> 8 ListBullet bullets (para 7 was empty, reused) ---
$psTexts = @(
  "Over 3.8 years of hands-on DevOps automation and CI/CD workflow management experience.",
  "DevOps Engineer at Cognizant, contributing to innovative solutions and operational efficiencies.",
  "Expertise in AWS infrastructure management using Terraform, Docker, and Kubernetes for deployments.",
  "Skilled in Jenkins, Git, GitLab, and Ansible for CI/CD pipeline optimization.",
  "Proven track record of automating deployment workflows, improving system reliability and performance.",
  "Experience with Prometheus and Grafana monitoring solutions for enhanced system visibility.",
  "Strong foundation in containerization, configuration management, and version control methodologies.",
  "Dedicated to high availability and operational excellence through automation strategies."
)
$p6 = $d.Paragraphs(6)
$p6.Style = "List Bullet"
$p6.Range.Text = $psTexts[0]

$p7 = $d.Paragraphs(7)
$p7.Style = "List Bullet"
$p7.Range.Text = $psTexts[1]

$insertAfter = $d.Paragraphs(7)
for ($i = 2; $i -lt $psTexts.Length; $i++) {
    $insertAfter.Range.InsertParagraphAfter()
    $newp = $d.Paragraphs(7 + ($i - 1))
    $newp.Style = "List Bullet"
    $newp.Range.Text = $psTexts[$i]
    $insertAfter = $newp
}

Write-Host "Step 10 done. Count: " $d.Paragraphs.Count
for ($i = 4; $i -le 16; $i++) {
    Write-Host $i ": [" $d.Paragraphs($i).Range.Text "]"
}

# --- 11. Insert new paragraph after "Contact:" line (para 2) ---
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(3)
$newPara.Range.Text = "Pursuing a career as a DevOps Engineer focused on enhancing Automation, Cloud Technologies, and Continuous Integration and delivery methodologies within innovative technology firms."

Write-Host "Step 11 done. Count: " $d.Paragraphs.Count
for ($i = 1; $i -le 6; $i++) {
    Write-Host $i ": [" $d.Paragraphs($i).Range.Text "]"
}
